# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the
# a8e4aae4-32ec-4080-b51d-80f27c8ba749 row (row 6) on both the
# "zh-cn" and "de-de" status sheets, recording the freshly generated
# handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-08 20:43:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-08 20:43:43"
